$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally held a "wide" table (years across columns, age bands down
# rows). This transposes it into a "long" table (years down rows, age bands
# across columns) headed by "Years", and left-aligns every cell.

# Clear the previous A1:G4 block entirely before laying out the new table.
$ws.Range("A1:G4").Clear()

$headers = @("Years", "        Under 18", "        18 to 24", "        Over 24 ")

$data = @(
    @(2015, 127787, 52973, 383948),
    @(2016, 120819, 50001, 379108),
    @(2017, 114529, 50992, 385475),
    @(2018, 111592, 48319, 392919),
    @(2019, 107069, 45629, 415017),
    @(2020, 106364, 45243, 428859)
)

# Header row
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Data rows
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}

# Left-align the whole new table (matches the new cellXfs entry applying
# horizontal="left").
$ws.Range("A1:D7").HorizontalAlignment = -4131  # xlLeft

# Column widths for the three new data columns (B:D); column A keeps its
# existing width of 15.
$ws.Columns.Item(2).ColumnWidth = 15.833333333333334  # ~16.71 chars
$ws.Columns.Item(3).ColumnWidth = 15.333333333333334  # ~16.14 chars
$ws.Columns.Item(4).ColumnWidth = 12.666666666666666  # ~13.57 chars

# Match the saved selection state.
$ws.Range("A1:D7").Select()
